$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# GL_Date (C2) - keep as text, not an Excel date
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "07/28/25"

# Vendor_Code (E2)
$ws.Range("E2").Value = "SLABRO"

# Invoice_Number (G2) - was text "110-S10112669.001", now numeric 860166732
$ws.Range("G2").Value = 860166732

# Invoice_Date (I2) - keep as text, not an Excel date
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "07/28/25"

# Invoice_Amount (J2)
$ws.Range("J2").Value = 81.53

# Remarks (R2)
$ws.Range("R2").Value = "Shop"

# Distribution_GL_Account (T2)
$ws.Range("T2").Value = 1200

# WO_Number (AA2) - clear value (was "2025"), cell remains as empty string
$ws.Range("AA2").Value = ""

# Item_Code (AB2) - remove value entirely (was "!Service Material")
$ws.Range("AB2").ClearContents()
